# 01/10/2023  Kevin Jackey  Completed Brand Revenue section
#
# Inserts a new validation row ("Include bots") into the Test Data Driver
# sheet, ahead of the existing "Saved Filters Create New Saved Filters
# button" row, widens column F to fit the new (longer) object-repo path,
# and updates the active selection to the new E11 cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 11 ("Saved Filters Create New
# Saved Filters button"); this shifts the old rows 11-17 down to 12-18.
$ws.Rows(11).Insert()

# Populate the new row 11 with the "Include bots" validation step.
$ws.Range("A11").Value = "validate"
$ws.Range("F11").Value = "Object Repository/Filters/Main Filter Window/Browsers and Devices/radio_IncludeBots"
$ws.Range("E11").Value = "Include bots"

# Column F now holds a longer object-repo path string, so widen it to fit.
$ws.Columns(6).ColumnWidth = 105.86

# Move / record the active selection on the new row.
$ws.Range("E11").Select()
